# The workbook already contains:
#   A1 = "all features numeric"        B1 = "non-numeric feature values"
#   A2 = "-"                           B2 = 0
#   A3 = 0                             B3 = 0
#   A4 = 0                             B4 = 0
#
# The only real content change in this revision is that the placeholder
# text "-" in A2 is removed (the cell becomes blank) while everything
# else on the sheet stays exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").ClearContents()
